$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Insert a new row above the current row 81 ("plot_end_time" / 2035), pushing
# everything below it down by one. The inserted row inherits the formatting
# of the row above (row 80), which already matches the desired style pattern
# (A: s30 string, B: s31 number, C/D: s69 blank, E: s30 blank).
$ws.Rows.Item(81).Insert()

# Populate the new row 81 with the new "plot_economics_start_time" parameter.
$ws.Cells.Item(81, 1).Value = "plot_economics_start_time"
$ws.Cells.Item(81, 2).Value = 2015

# Update the current selection/view to match where the edit was made.
$ws.Activate()
$ws.Range("E77").Select()
